# Add two new columns to the Movies bulk-upload template:
#   Q: is_available (boolean)
#   R: location (file path string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header formatting (bold, border, center/top alignment) from the
# existing "budget" header cell (P1) onto the two new header cells, then set
# their text. Copy first so the subsequent .Value assignment isn't clobbered
# by the paste.
$ws.Range("P1:P1").Copy($ws.Range("Q1"))
$ws.Range("P1:P1").Copy($ws.Range("R1"))

$ws.Range("Q1").Value = "is_available"
$ws.Range("R1").Value = "location"

# Row 2 - The Matrix
$ws.Range("Q2").Value = $true
$ws.Range("R2").Value = "/movies/the_matrix.mp4"

# Row 3 - Dangal
$ws.Range("Q3").Value = $true
$ws.Range("R3").Value = "/movies/dangal.mp4"

# Row 4 - Avengers: Endgame
$ws.Range("Q4").Value = $true
$ws.Range("R4").Value = "/movies/avengers_endgame.mp4"

# Row 5 - 3 Idiots
$ws.Range("Q5").Value = $true
$ws.Range("R5").Value = "/movies/3_idiots.mp4"

# Row 6 - Parasite (not available, no location on file)
$ws.Range("Q6").Value = $false
